# Auto-generated Excel COM-interop edit script
# Updates the cryptos list: refreshed prices / 1h volume deltas,
# plus a 3-way coin-row reorder (Monero/NEARProtocol and Stacks/Mantle/Maker).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.927.95'
$ws.Range('E2').Value = '  +2.96%  '
$ws.Range('D3').Value = '3.050.25'
$ws.Range('E3').Value = '  +2.35%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '525.07'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.68%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.10'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.49%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E8').Value = '  +4.82%  '
$ws.Range('E9').Value = '  +4.99%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.113'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.71%  '
$ws.Range('E11').Value = '  +5.16%  '
$ws.Range('E12').Value = '  +2.35%  '
$ws.Range('D13').Value = '3.573.60'
$ws.Range('E13').Value = '  +2.34%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.99'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +8.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000171'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +17.50%  '
$ws.Range('D16').Value = '57.895.34'
$ws.Range('E16').Value = '  +2.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.24'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +7.71%  '
$ws.Range('D18').Value = '3.049.37'
$ws.Range('E18').Value = '  +1.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.98'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.20'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '339.42'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.36%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.68'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.96%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.501'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +7.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.02'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.173'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +7.94%  '
$ws.Range('D27').Value = '0.0₃0974'
$ws.Range('E27').Value = '  +8.38%  '
$ws.Range('E28').Value = '  -0.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.01'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +8.65%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.35'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +8.86%  '
$ws.Range('E31').Value = '  +7.49%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.24'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.14'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.91%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.75'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.20%  '
$ws.Range('B35').Value = 'Monero'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '156.11'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.88%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.90'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.77%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.34'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.83%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.96'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +12.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0703'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.17%  '
$ws.Range('D40').Value = '3.086.31'
$ws.Range('E40').Value = '  +2.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '37.74'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.88'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +9.50%  '
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.344.53'
$ws.Range('E44').Value = '  +6.97%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.48'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.73%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.662'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.39%  '
$ws.Range('E47').Value = '  +3.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.03'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.72%  '
$ws.Range('E49').Value = '  +4.18%  '
$ws.Range('E50').Value = '  +5.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '20.17'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.48%  '
